# Apply the "2 outputs of one process" update to MasterTable.
# Net effect:
#  - Row 9 (O4) DESCRIPCION/USE_ACTUAL_MODEL changed to describe a new TANK1 input.
#  - A new O6 row is inserted right before the Y2 row (pushing Y2, X3, O6(old/-> O7) down).
#  - A new Y3 row (wrapped text, row height 45) is inserted right after the Y2 row.
#  - The last row's TAG/FEATURES_NAMES become O7 (new unique value) keeping the old
#    "proceso C" observed description/classification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update row 9 (O4) to describe the new TANK1 input ---------------
$ws.Range("C9").Value = "Variable de entrada al proceso tanque TANK1. No es una variable de decisión del optimizador. Variable Observada. Es un ejemplo de un flujo constante de entrada"
$ws.Range("F9").Value = "TANK1"

# --- 2. Insert a new row at 11 (duplicate of the old O6 "observed" row) -
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "O6"
$ws.Range("B11").Value = "O6"
$ws.Range("C11").Value = "Variable de entrada al proceso B. No es una variable de decisión del optimizador. Variable Observada"
$ws.Range("D11").Value = "Observed"
$ws.Range("E11").Value = "O"
$ws.Range("F11").Value = "MLB"

# --- 3. Insert a new row at 13 (new Y3 target row), after Y2 (now row 12)
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "Y3"
$ws.Range("B13").Value = "Y3"
$ws.Range("C13").Value = "Variable target del proceso B que finaliza el proceso"
$ws.Range("D13").Value = "Target"
$ws.Range("E13").Value = "T"
$ws.Range("F13").Value = "MLB"
$ws.Range("C13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 45

# --- 4. Final row (old O6/"proceso C" row, now row 15) becomes O7 -------
$ws.Range("A15").Value = "O7"
$ws.Range("B15").Value = "O7"

# --- 5. Sheet selection / dimension bookkeeping --------------------------
$ws.Range("A1:F15").Select()
